# Insert a new weekly record as row 734 (shifting the existing rows 734-759
# down to 735-760), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 734:759 down by one row (this also extends the sheet dimension
# to A1:R760 and keeps the date-format style on column D).
$ws.Rows("734:734").Insert()

# Populate the newly inserted row 734 with the new weekly record.
$ws.Cells.Item(734, 1).Value = 6
$ws.Cells.Item(734, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(734, 3).Value = "Metropolitana"
$ws.Cells.Item(734, 4).Value = 45075
$ws.Cells.Item(734, 5).Value = 13
$ws.Cells.Item(734, 6).Value = 100112012
$ws.Cells.Item(734, 7).Value = "Espinaca"
$ws.Cells.Item(734, 8).Value = "Sin especificar"
$ws.Cells.Item(734, 9).Value = "Primera"
$ws.Cells.Item(734, 10).Value = 340
$ws.Cells.Item(734, 11).Value = 7000
$ws.Cells.Item(734, 12).Value = 8000
$ws.Cells.Item(734, 13).Value = 7441
$ws.Cells.Item(734, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(734, 15).Value = "Región Metropolitana"
$ws.Cells.Item(734, 16).Value = 744
$ws.Cells.Item(734, 17).Value = 10
$ws.Cells.Item(734, 18).Value = "Hortaliza"
